$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.2274831784921
$ws.Range("C2").Value = 8.355915025247118
$ws.Range("D2").Value = 8.103667696083249
$ws.Range("E2").Value = 12.63965890866621
$ws.Range("F2").Value = 35.63091631070591
$ws.Range("J2").Value = 9.859911517634377
$ws.Range("L2").Value = 11.19252879243042
$ws.Range("N2").Value = 18.51592985720294
$ws.Range("O2").Value = 27.3655506043689
$ws.Range("B3").Value = 18.83105297692308
$ws.Range("C3").Value = 8.125909172112635
$ws.Range("D3").Value = 8.103615839439781
$ws.Range("E3").Value = 12.66858514225143
$ws.Range("F3").Value = 35.68512573369722
$ws.Range("J3").Value = 9.882424730450484
$ws.Range("L3").Value = 11.17846103245517
$ws.Range("N3").Value = 18.57153812877309
$ws.Range("O3").Value = 27.42774904111683
$ws.Range("B4").Value = 18.58641455618922
$ws.Range("C4").Value = 7.980213968521209
$ws.Range("D4").Value = 8.104414784910588
$ws.Range("E4").Value = 12.68773708461502
$ws.Range("F4").Value = 35.72695103253587
$ws.Range("J4").Value = 9.89700270388569
$ws.Range("L4").Value = 11.17119862931815
$ws.Range("N4").Value = 18.60755564934008
$ws.Range("O4").Value = 27.4721739788592
$ws.Range("B5").Value = 18.4865505647009
$ws.Range("C5").Value = 7.919775578427606
$ws.Range("D5").Value = 8.104950050786943
$ws.Range("E5").Value = 12.69589192634582
$ws.Range("F5").Value = 35.7461386551663
$ws.Range("J5").Value = 9.903133647265511
$ws.Range("L5").Value = 11.16858732881971
$ws.Range("N5").Value = 18.62270535606318
$ws.Range("O5").Value = 27.49184163534608
$ws.Range("B6").Value = 18.46996192502068
$ws.Range("C6").Value = 7.909677198482005
$ws.Range("D6").Value = 8.105051617931064
$ws.Range("E6").Value = 12.69726720059501
$ws.Range("F6").Value = 35.74945407147529
$ws.Range("J6").Value = 9.904163194405511
$ws.Range("L6").Value = 11.16817481990396
$ws.Range("N6").Value = 18.62524950981149
$ws.Range("O6").Value = 27.49520179084907
$ws.Range("B7").Value = 18.58506826740716
$ws.Range("C7").Value = 7.979403113445771
$ws.Range("D7").Value = 8.104421153664246
$ws.Range("E7").Value = 12.6878456448165
$ws.Range("F7").Value = 35.72720113068304
$ws.Range("J7").Value = 9.897084616736365
$ws.Range("L7").Value = 11.17116199960389
$ws.Range("N7").Value = 18.60775804991546
$ws.Range("O7").Value = 27.4724328957577
$ws.Range("B8").Value = 19.09113009239256
$ws.Range("C8").Value = 8.277570153677489
$ws.Range("D8").Value = 8.103477866897652
$ws.Range("E8").Value = 12.64934422573415
$ws.Range("F8").Value = 35.64783288696815
$ws.Range("J8").Value = 9.8675177271601
$ws.Range("L8").Value = 11.18739426166445
$ws.Range("N8").Value = 18.53471533193911
$ws.Range("O8").Value = 27.38570047220431
$ws.Range("B9").Value = 20.06796052008074
$ws.Range("C9").Value = 8.824495017270159
$ws.Range("D9").Value = 8.108184646511422
$ws.Range("E9").Value = 12.58486036538929
$ws.Range("F9").Value = 35.56011127419416
$ws.Range("J9").Value = 9.815502693231251
$ws.Range("L9").Value = 11.23002794669143
$ws.Range("N9").Value = 18.40629728028922
$ws.Range("O9").Value = 27.26524078556155
$ws.Range("B10").Value = 20.76874861954599
$ws.Range("C10").Value = 9.200517073293058
$ws.Range("D10").Value = 8.115590790846957
$ws.Range("E10").Value = 12.54417237425908
$ws.Range("F10").Value = 35.53723489104203
$ws.Range("J10").Value = 9.780891460653205
$ws.Range("L10").Value = 11.26778068874078
$ws.Range("N10").Value = 18.32091351637756
$ws.Range("O10").Value = 27.2071786534323
$ws.Range("B11").Value = 21.08244224020462
$ws.Range("C11").Value = 9.365476359181146
$ws.Range("D11").Value = 8.119805954673136
$ws.Range("E11").Value = 12.52710869349597
$ws.Range("F11").Value = 35.5358738226647
$ws.Range("J11").Value = 9.765921632414365
$ws.Range("L11").Value = 11.2863150046207
$ws.Range("N11").Value = 18.28400255283718
$ws.Range("O11").Value = 27.18740703136135
$ws.Range("B12").Value = 21.20038541216457
$ws.Range("C12").Value = 9.42702931914403
$ws.Range("D12").Value = 8.121522726769184
$ws.Range("E12").Value = 12.52085452089921
$ws.Range("F12").Value = 35.53665913752378
$ws.Range("J12").Value = 9.760363874769929
$ws.Range("L12").Value = 11.29352569144122
$ws.Range("N12").Value = 18.27030190274291
$ws.Range("O12").Value = 27.18087698082791
$ws.Range("B13").Value = 21.17502357314627
$ws.Range("C13").Value = 9.413813967951945
$ws.Range("D13").Value = 8.12114764409942
$ws.Range("E13").Value = 12.52219224761555
$ws.Range("F13").Value = 35.53643216589552
$ws.Range("J13").Value = 9.761555908091566
$ws.Range("L13").Value = 11.29196425253497
$ws.Range("N13").Value = 18.27324028894359
$ws.Range("O13").Value = 27.18224075180044
$ws.Range("B14").Value = 21.09216306636072
$ws.Range("C14").Value = 9.370558865849796
$ws.Range("D14").Value = 8.119944784308103
$ws.Range("E14").Value = 12.52659000310067
$ws.Range("F14").Value = 35.53591236591142
$ws.Range("J14").Value = 9.765462170486289
$ws.Range("L14").Value = 11.28690440183664
$ws.Range("N14").Value = 18.28286985096249
$ws.Range("O14").Value = 27.18685060963285
$ws.Range("B15").Value = 21.04129524655178
$ws.Range("C15").Value = 9.343943861340135
$ws.Range("D15").Value = 8.119223668148861
$ws.Range("E15").Value = 12.52931076294351
$ws.Range("F15").Value = 35.5357633498237
$ws.Range("J15").Value = 9.767869310022979
$ws.Range("L15").Value = 11.28383001386782
$ws.Range("N15").Value = 18.28880424857861
$ws.Range("O15").Value = 27.18979896795077
$ws.Range("B16").Value = 20.74813561609231
$ws.Range("C16").Value = 9.189610690609671
$ws.Range("D16").Value = 8.115332253692852
$ws.Range("E16").Value = 12.54531657982824
$ws.Range("F16").Value = 35.53750587875474
$ws.Range("J16").Value = 9.781885330273996
$ws.Range("L16").Value = 11.2665965086167
$ws.Range("N16").Value = 18.32336450641261
$ws.Range("O16").Value = 27.20860459694713
$ws.Range("B17").Value = 20.56690341627012
$ws.Range("C17").Value = 9.093343592224638
$ws.Range("D17").Value = 8.113160941916449
$ws.Range("E17").Value = 12.55550557964415
$ws.Range("F17").Value = 35.5408919178151
$ws.Range("J17").Value = 9.790681876461809
$ws.Range("L17").Value = 11.25637031919759
$ws.Range("N17").Value = 18.34505993020112
$ws.Range("O17").Value = 27.22184378925747
$ws.Range("B18").Value = 20.46218984696113
$ws.Range("C18").Value = 9.03740157465541
$ws.Range("D18").Value = 8.111991765148451
$ws.Range("E18").Value = 12.56150210238899
$ws.Range("F18").Value = 35.5436909425684
$ws.Range("J18").Value = 9.79581439007981
$ws.Range("L18").Value = 11.2506167520367
$ws.Range("N18").Value = 18.35772032661962
$ws.Range("O18").Value = 27.23008362216783
$ws.Range("B19").Value = 20.42665784089769
$ws.Range("C19").Value = 9.018363566856086
$ws.Range("D19").Value = 8.111609623433292
$ws.Range("E19").Value = 12.56355580553133
$ws.Range("F19").Value = 35.54478486958684
$ws.Range("J19").Value = 9.797564718021581
$ws.Range("L19").Value = 11.24869082631872
$ws.Range("N19").Value = 18.36203816435014
$ws.Range("O19").Value = 27.23298076283874
$ws.Range("B20").Value = 20.58624574825535
$ws.Range("C20").Value = 9.103650841749952
$ws.Range("D20").Value = 8.113383840351398
$ws.Range("E20").Value = 12.55440686230274
$ws.Range("F20").Value = 35.54044334193561
$ws.Range("J20").Value = 9.789737919968383
$ws.Range("L20").Value = 11.2574456641111
$ws.Range("N20").Value = 18.3427316099923
$ws.Range("O20").Value = 27.22036975327385
$ws.Range("B21").Value = 21.1165250221495
$ws.Range("C21").Value = 9.383289007766479
$ws.Range("D21").Value = 8.120294829737178
$ws.Range("E21").Value = 12.52529264803801
$ws.Range("F21").Value = 35.53602974718201
$ws.Range("J21").Value = 9.764311797834219
$ws.Range("L21").Value = 11.28838541646435
$ws.Range("N21").Value = 18.2800339116669
$ws.Range("O21").Value = 27.1854705948886
$ws.Range("B22").Value = 21.45811613888675
$ws.Range("C22").Value = 9.560709332482565
$ws.Range("D22").Value = 8.125513860895978
$ws.Range("E22").Value = 12.50747395768505
$ws.Range("F22").Value = 35.54072615581273
$ws.Range("J22").Value = 9.748341093843889
$ws.Range("L22").Value = 11.30972456494288
$ws.Range("N22").Value = 18.24066994858188
$ws.Range("O22").Value = 27.16824093887672
$ws.Range("B23").Value = 21.27629326606196
$ws.Range("C23").Value = 9.466516597718204
$ws.Range("D23").Value = 8.122664476235601
$ws.Range("E23").Value = 12.51687362026422
$ws.Range("F23").Value = 35.53752617812464
$ws.Range("J23").Value = 9.756805930751389
$ws.Range("L23").Value = 11.29823433327501
$ws.Range("N23").Value = 18.26153196710729
$ws.Range("O23").Value = 27.17692569308173
$ws.Range("B24").Value = 20.57750269335241
$ws.Range("C24").Value = 9.098992790265891
$ws.Range("D24").Value = 8.113282821389955
$ws.Range("E24").Value = 12.55490316025446
$ws.Range("F24").Value = 35.54064348806587
$ws.Range("J24").Value = 9.790164448227234
$ws.Range("L24").Value = 11.25695910916954
$ws.Range("N24").Value = 18.34378365970842
$ws.Range("O24").Value = 27.22103420750038
$ws.Range("B25").Value = 19.80615498248556
$ws.Range("C25").Value = 8.680885141230418
$ws.Range("D25").Value = 8.106214270931881
$ws.Range("E25").Value = 12.60112839112668
$ws.Range("F25").Value = 35.57655051307744
$ws.Range("J25").Value = 9.828938867120781
$ws.Range("L25").Value = 11.21735331956647
$ws.Range("N25").Value = 18.43945853815873
$ws.Range("O25").Value = 27.29249494016692
